# Revert "Powerpoint writer: consolidate text run nodes."
#
# The caption textbox on slide 1 currently stores "The picture first" as
# three runs: "The ", "picture ", "first". The target (pre-consolidation)
# layout splits the leading/trailing spaces of the first two words into
# their own runs: "The", " ", "picture", " ", "first".
#
# Setting .Text on a Characters() sub-range that lines up with the desired
# run boundaries causes the writer to split the run there (even when the
# replacement text equals the original), without touching neighboring runs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Full text is "The picture first" (chars 1-18):
#   1-3  "The"       4  " "
#   5-11 "picture"  12  " "
#  13-17 "first"   (left untouched)

$tr.Characters(1, 3).Text = "The"
$tr.Characters(4, 1).Text = " "
$tr.Characters(5, 7).Text = "picture"
$tr.Characters(12, 1).Text = " "
